# verigen : lua 리스트 추가 ("verigen_description" 시트 추가)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. summary sheet (sheet1): register the new "verigen_description"
#    function in the summary table (row 4).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Range("A4").Value = "verigen_description"
$summary.Range("B4").Value = "function"
$summary.Range("B4").HorizontalAlignment = -4108   # xlCenter, matches B2/B3
$summary.Range("C4").Value = "Set verigen source code description"
$summary.Range("A4").Select()

# ---------------------------------------------------------------------
# 2. Build the new "verigen_description" sheet right after "vfunction"
#    by duplicating "vfunction" (same look & feel: columns / styles)
#    and then overwriting its content.
# ---------------------------------------------------------------------
$vfunction = $wb.Worksheets.Item(4)
$vfunction.Copy($null, $vfunction)

$descSheet = $wb.Worksheets.Item(5)
$descSheet.Name = "verigen_description"

$descSheet.Range("B2").Value = "function verigen_description(desc)"
$descSheet.Range("B4").Value = "Set verigen source's description in Lua file"
$descSheet.Rows.Item(4).AutoFit()
$descSheet.Range("A5").Value = "desc"
$descSheet.Range("B5").Value = "description of verigen source file"
$descSheet.Rows.Item(6).Delete()

# ---------------------------------------------------------------------
# 3. Tidy up selections left over from the copy/edit operations.
# ---------------------------------------------------------------------
$vfunction.Range("A1:B5").Select()
$descSheet.Range("A1").Select()

Write-Host "verigen_description sheet added"
